$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '31.262.02'
$ws.Range('E2').Value = '  +2.98%  '

# Row 3
$ws.Range('D3').Value = '1.988.21'
$ws.Range('E3').Value = '  +6.24%  '

# Row 4
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.17%  '

# Row 5
$ws.Range('D5').Value = '0.8105'
$ws.Range('E5').Value = '  +72.05%  '

# Row 6
$ws.Range('D6').Value = '253.80'
$ws.Range('E6').Value = '  +4.06%  '

# Row 7
$ws.Range('D7').Value = '0.9978'
$ws.Range('E7').Value = '  -0.24%  '

# Row 8
$ws.Range('D8').Value = '0.3439'

# Row 9
$ws.Range('D9').Value = '25.63'
$ws.Range('E9').Value = '  +16.57%  '

# Row 10
$ws.Range('D10').Value = '0.06966'
$ws.Range('E10').Value = '  +7.91%  '

# Row 11
$ws.Range('D11').Value = '0.8421'
$ws.Range('E11').Value = '  +16.08%  '

# Row 12
$ws.Range('D12').Value = '0.08113'
$ws.Range('E12').Value = '  +4.34%  '

# Row 13
$ws.Range('D13').Value = '1.987.45'
$ws.Range('E13').Value = '  +6.17%  '

# Row 14
$ws.Range('D14').Value = '100.75'
$ws.Range('E14').Value = '  +4.73%  '

# Row 15
$ws.Range('D15').Value = '5.512'
$ws.Range('E15').Value = '  +7.42%  '

# Row 16
$ws.Range('D16').Value = '272.49'
$ws.Range('E16').Value = '  -2.40%  '

# Row 17
$ws.Range('D17').Value = '31.251.32'
$ws.Range('E17').Value = '  +2.99%  '

# Row 18
$ws.Range('D18').Value = '13.97'
$ws.Range('E18').Value = '  +7.35%  '

# Row 19
$ws.Range('D19').Value = '0.000007927'
$ws.Range('E19').Value = '  +5.71%  '

# Row 20
$ws.Range('D20').Value = '5.814'
$ws.Range('E20').Value = '  +10.80%  '

# Row 21
$ws.Range('D21').Value = '2.247.36'
$ws.Range('E21').Value = '  +6.34%  '

# Row 22
$ws.Range('D22').Value = '0.9974'
$ws.Range('E22').Value = '  -0.28%  '

# Row 23
$ws.Range('D23').Value = '0.9984'
$ws.Range('E23').Value = '  -0.19%  '

# Row 24
$ws.Range('D24').Value = '6.942'
$ws.Range('E24').Value = '  +11.41%  '

# Row 25
$ws.Range('D25').Value = '9.763'
$ws.Range('E25').Value = '  +7.78%  '

# Row 26: 'Monero' -> 'Stellar'
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = '0.1502'
$ws.Range('E26').Value = '  +55.99%  '

# Row 27: 'Stellar' -> 'Monero'
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '163.95'
$ws.Range('E27').Value = '  +0.16%  '

# Row 28
$ws.Range('D28').Value = '20.07'
$ws.Range('E28').Value = '  +7.36%  '

# Row 29
$ws.Range('D29').Value = '2.183'
$ws.Range('E29').Value = '  +16.33%  '

# Row 30
$ws.Range('D30').Value = '1.566'
$ws.Range('E30').Value = '  +5.45%  '

# Row 31: 'Toncoin' -> 'Filecoin'
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '4.566'
$ws.Range('E31').Value = '  +8.26%  '

# Row 32: 'Filecoin' -> 'Toncoin'
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '1.346'
$ws.Range('E32').Value = '  +1.93%  '

# Row 33
$ws.Range('E33').Value = '  +4.98%  '

# Row 34
$ws.Range('D34').Value = '0.05167'
$ws.Range('E34').Value = '  +7.57%  '

# Row 35
$ws.Range('D35').Value = '1.214'
$ws.Range('E35').Value = '  +8.37%  '

# Row 36
$ws.Range('D36').Value = '0.7576'
$ws.Range('E36').Value = '  +10.13%  '

# Row 37
$ws.Range('D37').Value = '2.765'
$ws.Range('E37').Value = '  +2.07%  '

# Row 38
$ws.Range('E38').Value = '  +6.21%  '

# Row 39
$ws.Range('D39').Value = '2.912'
$ws.Range('E39').Value = '  +3.65%  '

# Row 40
$ws.Range('D40').Value = '6.600'
$ws.Range('E40').Value = '  +6.06%  '

# Row 41: 'Aave' -> 'TheSandbox'
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.4698'
$ws.Range('E41').Value = '  +11.28%  '

# Row 42: 'TheSandbox' -> 'Aave'
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '78.02'
$ws.Range('E42').Value = '  +5.12%  '

# Row 43
$ws.Range('D43').Value = '2.074'
$ws.Range('E43').Value = '  +7.41%  '

# Row 44
$ws.Range('D44').Value = '0.8516'
$ws.Range('E44').Value = '  +3.60%  '

# Row 45
$ws.Range('D45').Value = '104.68'
$ws.Range('E45').Value = '  +3.71%  '

# Row 46
$ws.Range('D46').Value = '0.9977'
$ws.Range('E46').Value = '  -0.18%  '

# Row 47
$ws.Range('D47').Value = '9.930'
$ws.Range('E47').Value = '  +3.31%  '

# Row 48
$ws.Range('D48').Value = '7.503'
$ws.Range('E48').Value = '  +7.85%  '

# Row 49
$ws.Range('D49').Value = '0.4303'
$ws.Range('E49').Value = '  +10.02%  '

# Row 50
$ws.Range('D50').Value = '36.69'
$ws.Range('E50').Value = '  +4.09%  '

# Row 51
$ws.Range('D51').Value = '0.1188'
$ws.Range('E51').Value = '  +12.11%  '
